# Generate Report for Handoff
# Adds a new handoff record (cbf06a04-...) as row 3 on each of the three
# worksheets (Overview, zh-cn, de-de), growing each Excel Table by one row.

$wb = $excel.ActiveWorkbook

# ---- shared literal values -------------------------------------------------
$newFileName   = 'cbf06a04-2fdc-4c4f-9f41-80508137ede6oooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$newFileDisp   = 'e2e\cbf06a04-2fdc-4c4f-9f41-80508137ede6oooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$newFileUrl    = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6123740b3a0279c56e2aeaaf5e9db197a62b643b/e2e/cbf06a04-2fdc-4c4f-9f41-80508137ede6oooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'

$extMd         = '.md'
$readyStatus   = 'Ready for handoff'
$overviewDate  = '2016-11-09 01:59:31'

$e2e           = 'e2e'
$htToken       = 'ht'
$falseStr      = 'False'
$trueStr       = 'True'
$epochDate     = '0001-01-01 00:00:00'

$zhXliff       = 'cbf06a04-2fdc-4c4f-9f41-80508137ede6ooooooooooooooooooooooooooooooooooooo.fa052df3b5c4b52e7ad1813ec2b4a9da4dcc0311.zh-cn.xlf'
$zhXliffDate   = '2016-11-09 01:59:17'

$deXliff       = 'cbf06a04-2fdc-4c4f-9f41-80508137ede6ooooooooooooooooooooooooooooooooooooo.fa052df3b5c4b52e7ad1813ec2b4a9da4dcc0311.de-de.xlf'
$deXliffDate   = '2016-11-09 01:59:31'

$dateFormat    = 'yyyy-mm-dd HH:mm:ss'

# =============================================================================
# Sheet 1: "Overview"
# =============================================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newFileName
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newFileUrl, "", "", $newFileDisp) | Out-Null
$wsOverview.Range("C3").Value = $extMd
$wsOverview.Range("E3").Value = $readyStatus
$wsOverview.Range("F3").Value = $readyStatus
$wsOverview.Range("G3").Value = $overviewDate
$wsOverview.Range("G3").NumberFormat = $dateFormat

# ---- column widths grow because the "Ready for handoff" text is longer ----
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null

# =============================================================================
# Sheet 2: "zh-cn"
# =============================================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newFileUrl, "", "", $newFileName) | Out-Null
$wsZhCn.Range("B3").Value = $extMd
$wsZhCn.Range("C3").Value = $readyStatus
$wsZhCn.Range("D3").Value = $e2e
$wsZhCn.Range("E3").Value = $htToken
$wsZhCn.Range("F3").Value = $falseStr
$wsZhCn.Range("G3").Value = $zhXliff
$wsZhCn.Range("H3").Value = $zhXliffDate
$wsZhCn.Range("H3").NumberFormat = $dateFormat
$wsZhCn.Range("K3").Value = $epochDate
$wsZhCn.Range("K3").NumberFormat = $dateFormat
$wsZhCn.Range("M3").Value = $trueStr
$wsZhCn.Range("O3").Value = $falseStr

$wsZhCn.Columns.Item(3).AutoFit() | Out-Null

# =============================================================================
# Sheet 3: "de-de"
# =============================================================================
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newFileUrl, "", "", $newFileName) | Out-Null
$wsDeDe.Range("B3").Value = $extMd
$wsDeDe.Range("C3").Value = $readyStatus
$wsDeDe.Range("D3").Value = $e2e
$wsDeDe.Range("E3").Value = $htToken
$wsDeDe.Range("F3").Value = $falseStr
$wsDeDe.Range("G3").Value = $deXliff
$wsDeDe.Range("H3").Value = $deXliffDate
$wsDeDe.Range("H3").NumberFormat = $dateFormat
$wsDeDe.Range("K3").Value = $epochDate
$wsDeDe.Range("K3").NumberFormat = $dateFormat
$wsDeDe.Range("M3").Value = $trueStr
$wsDeDe.Range("O3").Value = $falseStr

$wsDeDe.Columns.Item(3).AutoFit() | Out-Null

Write-Output "Added handoff row for cbf06a04 on Overview, zh-cn and de-de sheets."
